$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: 4x100m - M results
$ws.Range("B8").Value = "CAN"
$ws.Range("C8").Value = "RSA"

# Row 7: 400m - F results
$ws.Range("B7").Value = "Marileidy Paulino"
$ws.Range("C7").Value = "Salwa Eid Naser"
$ws.Range("D7").Value = "Natalia Kaczmarek"

# Row 8 (3rd place) - reuses existing "GBR" string
$ws.Range("D8").Value = "GBR"

$ws.Range("D8").Select()
